$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$values = @(
  "Sát Phá Tham hội chiếu tại Tử Tức",
  "Tử Phủ Vũ Tướng hội chiếu tại Tử Tức",
  "Cơ Nguyệt Đồng Lương hội chiếu tại Tử Tức",
  "Cự Nhật hội chiếu tại Tử Tức",
  "Kình Đà hội chiếu tại Tử Tức",
  "Xương Khúc hội chiếu tại Tử Tức",
  "Hoả Linh hội chiếu tại Tử Tức",
  "Không Kiếp hội chiếu tại Tử Tức",
  "Quang Quý hội chiếu tại Tử Tức",
  "Tả Hữu hội chiếu tại Tử Tức",
  "Song Hao hội chiếu tại Tử Tức",
  "Tang Hổ hội chiếu tại Tử Tức",
  "Khốc Hư hội chiếu tại Tử Tức",
  "Hình Riêu hội chiếu tại Tử Tức",
  "Thai Toạ hội chiếu tại Tử Tức",
  "Đào Hồng hội chiếu tại Tử Tức",
  "Ấn Phù hội chiếu tại Tử Tức",
  "Song Hao Quyền Lộc Kiếp Hoả hội chiếu tại Tử Tức",
  "Tử Phủ Vũ Tướng Xương Khúc Khôi Việt Tả Hữu Khoa Quyền Lộc Long hội chiếu tại Tử Tức",
  "Tử Phủ Vũ Tướng Tả Hữu Khoa Quyền Lộc Long Phượng hội chiếu tại Tử Tức",
  "Tử Khúc Phá Dương Đà hội chiếu tại Tử Tức",
  "Cơ Nguyệt Đồng Lương gặp Xương Khúc Tả Hữu hội chiếu tại Tử Tức",
  "Cơ Nguyệt Đồng Lương Khoa Tả Hữu Quang Quý Quan Phúc hội chiếu tại Tử Tức",
  "Sát Quyền hội chiếu tại Tử Tức",
  "Lộc Mã hội chiếu tại Tử Tức",
  "Kiếp Hư Hao Quyền hội chiếu tại Tử Tức",
  "Tuế Hổ Phù Xương Khúc hội chiếu tại Tử Tức",
  "Xương Khúc Tấu Long Phượng hội chiếu tại Tử Tức",
  "Đào Hồng Riêu Tấu Cơ Vũ hội chiếu tại Tử Tức",
  "Binh Hình Tướng Ấn hội chiếu tại Tử Tức",
  "Hổ Tấu hội chiếu tại Tử Tức",
  "Hình Riêu Y hội chiếu tại Tử Tức",
  "Mã Hỏa Linh hội chiếu tại Tử Tức",
  "Thai Tọa Hồng Đào hội chiếu tại Tử Tức",
  "Tả Hữu Không Kiếp hội chiếu tại Tử Tức",
  "Tả Hữu Binh Tướng hội chiếu tại Tử Tức",
  "Đào Quyền hội chiếu tại Tử Tức",
  "Đào Hồng Tả Cái Hữu hội chiếu tại Tử Tức",
  "Quan Phúc Quang Tấu hội chiếu tại Tử Tức",
  "Đào Hồng Xương Khúc hội chiếu tại Tử Tức",
  "Đào Tử Phủ hội chiếu tại Tử Tức",
  "Hổ Kình Sát hội chiếu tại Tử Tức",
  "Hổ Tang Kiếp hội chiếu tại Tử Tức",
  "Hổ Tang Không Kiếp hội chiếu tại Tử Tức",
  "Âm Dương Lương hội chiếu tại Tử Tức",
  "Cơ Lương Gia Hội hội chiếu tại Tử Tức",
  "Nhật Chiếu Lôi Môn hội chiếu tại Tử Tức",
  "Tả Hữu Xương Khúc hội chiếu tại Tử Tức",
  "Tham Linh Triều Viên hội chiếu tại Tử Tức",
  "Cự Hỏa Linh hội chiếu tại Tử Tức",
  "Hồng Đào Kỵ hội chiếu tại Tử Tức",
  "Cự Đồng Hình hội chiếu tại Tử Tức",
  "Lương Phá hội chiếu tại Tử Tức",
  "Khôi Việt hội chiếu tại Tử Tức",
  "Tham Vũ Hỏa hội chiếu tại Tử Tức",
  "Cự Kỵ hội chiếu tại Tử Tức",
  "Kình Đà Hỏa hội chiếu tại Tử Tức",
  "Kình Đà Không Kiếp hội chiếu tại Tử Tức",
  "Đào Hồng Hỷ hội chiếu tại Tử Tức",
  "Đào Hồng Kiếp Sát hội chiếu tại Tử Tức",
  "Đào Hồng Kỵ hội chiếu tại Tử Tức",
  "Đào Hồng Tả Phù Hữu Bật hội chiếu tại Tử Tức",
  "Đào Hồng Tả Phù Hữu Bật Khoa Quyền Lộc hội chiếu tại Tử Tức",
  "Đào Hồng Tả Phù Hữu Bật Khoa Quyền Lộc Long Trì Phượng Các hội chiếu tại Tử Tức",
  "Lưu Hà Kiếp Sát hội chiếu tại Tử Tức",
  "Phục Không Kiếp hội chiếu tại Tử Tức",
  "Xương Khúc Khôi Việt hội chiếu tại Tử Tức",
  "Khoa Quyền hội chiếu tại Tử Tức",
  "Lộc Quyền hội chiếu tại Tử Tức",
  "Tướng Binh Đào Hồng hội chiếu tại Tử Tức",
  "Lương Khốc Tuế hội chiếu tại Tử Tức",
  "Tướng Binh hội chiếu tại Tử Tức",
  "Cự Tang hội chiếu tại Tử Tức",
  "Cự Hỏa hội chiếu tại Tử Tức",
  "Phá Hình Kỵ hội chiếu tại Tử Tức",
  "Tang Trực Tuế hội chiếu tại Tử Tức",
  "Tuế Xương Khúc hội chiếu tại Tử Tức",
  "Tử Vi Tả Hữu hội chiếu tại Tử Tức",
  "Tử Vi Tang Tả Hữu hội chiếu tại Tử Tức"
)
$startRow = 4225
for ($i = 0; $i -lt $values.Count; $i++) {
  $r = $startRow + $i
  $ws.Cells.Item($r, 1).Value = $values[$i]
  $ws.Cells.Item($r, 2).Value = $values[$i]
}